$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new simulation runs ("Holden" and "Rizzie Spiral") were inserted as rows 4-5,
# pushing the previously-existing runs (rows 4-29) down to rows 6-31, and the
# simulation was rerun producing new numeric results throughout. "Thomas Hex" was
# also renamed to "Matthies Hex".

# --- New row index (column A) for the two brand-new trailing rows ---
$ws.Cells.Item(30,1).Value = 28
$ws.Cells.Item(31,1).Value = 29

# Copy the bold/bordered style used by column A onto the two new rows
$ws.Range("A29").Copy() | Out-Null
$ws.Range("A30:A31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Column B (run labels) for rows 4-31 ---
$ws.Cells.Item(4,2).Value = 'Holden'
$ws.Cells.Item(5,2).Value = 'Rizzie Spiral'
$ws.Cells.Item(6,2).Value = 'RotRing OmegaMax-90'
$ws.Cells.Item(7,2).Value = 'Equal Angle'
$ws.Cells.Item(8,2).Value = 'Tilt Rotate'
$ws.Cells.Item(9,2).Value = 'CLR'
$ws.Cells.Item(10,2).Value = 'Rizzie Hex'
$ws.Cells.Item(11,2).Value = 'Matthies Hex'
$ws.Cells.Item(12,2).Value = 'Tilt Rotate_Partial'
$ws.Cells.Item(13,2).Value = 'RotRing OmegaMax-60'
$ws.Cells.Item(14,2).Value = 'Equal Angle_Partial'
$ws.Cells.Item(15,2).Value = 'Rizzie Hex_Partial'
$ws.Cells.Item(16,2).Value = 'ND Single'
$ws.Cells.Item(17,2).Value = 'RD Single'
$ws.Cells.Item(18,2).Value = 'TD Single'
$ws.Cells.Item(19,2).Value = 'Morris Single'
$ws.Cells.Item(20,2).Value = 'Ring Perpendicular to ND'
$ws.Cells.Item(21,2).Value = 'Ring Perpendicular to RD'
$ws.Cells.Item(22,2).Value = 'Ring Perpendicular to TD'
$ws.Cells.Item(23,2).Value = 'OffsetFTD'
$ws.Cells.Item(24,2).Value = 'OffsetATD'
$ws.Cells.Item(25,2).Value = 'OffsetF45'
$ws.Cells.Item(26,2).Value = 'OffsetA45'
$ws.Cells.Item(27,2).Value = 'OffsetFRD'
$ws.Cells.Item(28,2).Value = 'OffsetARD'
$ws.Cells.Item(29,2).Value = 'Gaussian Quadrature'
$ws.Cells.Item(30,2).Value = 'Michael-CCHex'
$ws.Cells.Item(31,2).Value = 'Michael-SNHex'

# --- Data columns C:T for rows 4-31 (rerun simulation results) ---
$ws.Cells.Item(4,3).Value = 1.006560340931292
$ws.Cells.Item(4,4).Value = 0.9270970518834148
$ws.Cells.Item(4,5).Value = 1.120621709077693
$ws.Cells.Item(4,6).Value = 1.031402538240342
$ws.Cells.Item(4,7).Value = 1.041625661685543
$ws.Cells.Item(4,8).Value = 0.8843346339685189
$ws.Cells.Item(4,9).Value = 0.8843346339685189
$ws.Cells.Item(4,10).Value = 1.006560340931292
$ws.Cells.Item(4,11).Value = 1.006560340931292
$ws.Cells.Item(4,12).Value = 1.041625661685543
$ws.Cells.Item(4,13).Value = 0.9629801478270308
$ws.Cells.Item(4,14).Value = 0.9629801478270308
$ws.Cells.Item(4,15).Value = 0.9510191158458254
$ws.Cells.Item(4,16).Value = 0.9775068788617846
$ws.Cells.Item(4,17).Value = 0.9775068788617846
$ws.Cells.Item(4,18).Value = 0.9847702443791615
$ws.Cells.Item(4,19).Value = 0.9847702443791615
$ws.Cells.Item(4,20).Value = 1.001940322631134
$ws.Cells.Item(5,3).Value = 0.956850357467216
$ws.Cells.Item(5,4).Value = 0.9505915255347241
$ws.Cells.Item(5,5).Value = 1.142382644922473
$ws.Cells.Item(5,6).Value = 1.017268461113676
$ws.Cells.Item(5,7).Value = 1.042131864647791
$ws.Cells.Item(5,8).Value = 0.9608000721303295
$ws.Cells.Item(5,9).Value = 0.9608000721303295
$ws.Cells.Item(5,10).Value = 0.956850357467216
$ws.Cells.Item(5,11).Value = 0.956850357467216
$ws.Cells.Item(5,12).Value = 1.042131864647791
$ws.Cells.Item(5,13).Value = 1.00146596838906
$ws.Cells.Item(5,14).Value = 1.00146596838906
$ws.Cells.Item(5,15).Value = 0.9845078207709482
$ws.Cells.Item(5,16).Value = 0.9865940980817788
$ws.Cells.Item(5,17).Value = 0.9865940980817788
$ws.Cells.Item(5,18).Value = 0.9791581629281381
$ws.Cells.Item(5,19).Value = 0.9791581629281381
$ws.Cells.Item(5,20).Value = 1.011670820969368
$ws.Cells.Item(6,3).Value = 1.036954644247129
$ws.Cells.Item(6,4).Value = 0.9713806539614657
$ws.Cells.Item(6,5).Value = 1.008752713212554
$ws.Cells.Item(6,6).Value = 1.013876379666414
$ws.Cells.Item(6,7).Value = 1.004379151895076
$ws.Cells.Item(6,8).Value = 0.9375382428229071
$ws.Cells.Item(6,9).Value = 0.9375382428229071
$ws.Cells.Item(6,10).Value = 1.036954644247129
$ws.Cells.Item(6,11).Value = 1.036954644247129
$ws.Cells.Item(6,12).Value = 1.004379151895076
$ws.Cells.Item(6,13).Value = 0.9709586973589918
$ws.Cells.Item(6,14).Value = 0.9709586973589918
$ws.Cells.Item(6,15).Value = 0.9710993495598164
$ws.Cells.Item(6,16).Value = 0.992957346321704
$ws.Cells.Item(6,17).Value = 0.992957346321704
$ws.Cells.Item(6,18).Value = 1.00395667080306
$ws.Cells.Item(6,19).Value = 1.00395667080306
$ws.Cells.Item(6,20).Value = 0.9954802976342577
$ws.Cells.Item(7,3).Value = 1.008810256138329
$ws.Cells.Item(7,4).Value = 0.9130316535518725
$ws.Cells.Item(7,5).Value = 1.142954550252163
$ws.Cells.Item(7,6).Value = 1.037466948242076
$ws.Cells.Item(7,7).Value = 1.04925982072046
$ws.Cells.Item(7,8).Value = 0.8618115352017307
$ws.Cells.Item(7,9).Value = 0.8618115352017307
$ws.Cells.Item(7,10).Value = 1.008810256138329
$ws.Cells.Item(7,11).Value = 1.008810256138329
$ws.Cells.Item(7,12).Value = 1.04925982072046
$ws.Cells.Item(7,13).Value = 0.9555356779610953
$ws.Cells.Item(7,14).Value = 0.9555356779610953
$ws.Cells.Item(7,15).Value = 0.9413676698246878
$ws.Cells.Item(7,16).Value = 0.9732938706868398
$ws.Cells.Item(7,17).Value = 0.9732938706868399
$ws.Cells.Item(7,18).Value = 0.9821729670497121
$ws.Cells.Item(7,19).Value = 0.9821729670497121
$ws.Cells.Item(7,20).Value = 1.002222460684438
$ws.Cells.Item(8,3).Value = 0.9839825738908773
$ws.Cells.Item(8,4).Value = 0.7100076925596597
$ws.Cells.Item(8,5).Value = 1.533300751186182
$ws.Cells.Item(8,6).Value = 1.121762140148446
$ws.Cells.Item(8,7).Value = 1.178461508073961
$ws.Cells.Item(8,8).Value = 0.5680991623450676
$ws.Cells.Item(8,9).Value = 0.5680991623450676
$ws.Cells.Item(8,10).Value = 0.9839825738908773
$ws.Cells.Item(8,11).Value = 0.9839825738908773
$ws.Cells.Item(8,12).Value = 1.178461508073961
$ws.Cells.Item(8,13).Value = 0.8732803352095144
$ws.Cells.Item(8,14).Value = 0.8732803352095144
$ws.Cells.Item(8,15).Value = 0.8188561209928963
$ws.Cells.Item(8,16).Value = 0.9101810814366353
$ws.Cells.Item(8,17).Value = 0.9101810814366353
$ws.Cells.Item(8,18).Value = 0.9286314545501958
$ws.Cells.Item(8,19).Value = 0.9286314545501958
$ws.Cells.Item(8,20).Value = 1.015935638034032
$ws.Cells.Item(9,3).Value = 1.007762072131477
$ws.Cells.Item(9,4).Value = 0.9906628297829714
$ws.Cells.Item(9,5).Value = 1.00735312536778
$ws.Cells.Item(9,6).Value = 1.004397526991071
$ws.Cells.Item(9,7).Value = 1.003012940847726
$ws.Cells.Item(9,8).Value = 0.9813991357291133
$ws.Cells.Item(9,9).Value = 0.9813991357291133
$ws.Cells.Item(9,10).Value = 1.007762072131477
$ws.Cells.Item(9,11).Value = 1.007762072131477
$ws.Cells.Item(9,12).Value = 1.003012940847726
$ws.Cells.Item(9,13).Value = 0.9922060382884195
$ws.Cells.Item(9,14).Value = 0.9922060382884195
$ws.Cells.Item(9,15).Value = 0.9916916354532702
$ws.Cells.Item(9,16).Value = 0.997391382902772
$ws.Cells.Item(9,17).Value = 0.997391382902772
$ws.Cells.Item(9,18).Value = 0.9999840552099482
$ws.Cells.Item(9,19).Value = 0.9999840552099482
$ws.Cells.Item(9,20).Value = 0.9990979384750229
$ws.Cells.Item(10,3).Value = 0.9991815967586943
$ws.Cells.Item(10,4).Value = 0.99931143537058
$ws.Cells.Item(10,5).Value = 1.002033769140755
$ws.Cells.Item(10,6).Value = 1.000277653754377
$ws.Cells.Item(10,7).Value = 1.000728706520544
$ws.Cells.Item(10,8).Value = 0.9992038155274823
$ws.Cells.Item(10,9).Value = 0.9992038155274823
$ws.Cells.Item(10,10).Value = 0.9991815967586943
$ws.Cells.Item(10,11).Value = 0.9991815967586943
$ws.Cells.Item(10,12).Value = 1.000728706520544
$ws.Cells.Item(10,13).Value = 0.9999662610240133
$ws.Cells.Item(10,14).Value = 0.9999662610240133
$ws.Cells.Item(10,15).Value = 0.9997479858062022
$ws.Cells.Item(10,16).Value = 0.999704706268907
$ws.Cells.Item(10,17).Value = 0.999704706268907
$ws.Cells.Item(10,18).Value = 0.9995739288913539
$ws.Cells.Item(10,19).Value = 0.9995739288913539
$ws.Cells.Item(10,20).Value = 1.000122829512072
$ws.Cells.Item(11,3).Value = 1.014134673881566
$ws.Cells.Item(11,4).Value = 0.9834711376313343
$ws.Cells.Item(11,5).Value = 1.012612920375124
$ws.Cells.Item(11,6).Value = 1.007792802506031
$ws.Cells.Item(11,7).Value = 1.005182977435622
$ws.Cells.Item(11,8).Value = 0.9669519565561844
$ws.Cells.Item(11,9).Value = 0.9669519565561844
$ws.Cells.Item(11,10).Value = 1.014134673881566
$ws.Cells.Item(11,11).Value = 1.014134673881566
$ws.Cells.Item(11,12).Value = 1.005182977435622
$ws.Cells.Item(11,13).Value = 0.9860674669959031
$ws.Cells.Item(11,14).Value = 0.9860674669959031
$ws.Cells.Item(11,15).Value = 0.9852020238743803
$ws.Cells.Item(11,16).Value = 0.9954232026244574
$ws.Cells.Item(11,17).Value = 0.9954232026244574
$ws.Cells.Item(11,18).Value = 1.000101070438735
$ws.Cells.Item(11,19).Value = 1.000101070438735
$ws.Cells.Item(11,20).Value = 0.9983577447309768
$ws.Cells.Item(12,3).Value = 0.9814795055814179
$ws.Cells.Item(12,4).Value = 0.7057314585161016
$ws.Cells.Item(12,5).Value = 1.544159252961416
$ws.Cells.Item(12,6).Value = 1.123344991014803
$ws.Cells.Item(12,7).Value = 1.181725805954586
$ws.Cells.Item(12,8).Value = 0.5637216591357709
$ws.Cells.Item(12,9).Value = 0.5637216591357709
$ws.Cells.Item(12,10).Value = 0.9814795055814179
$ws.Cells.Item(12,11).Value = 0.9814795055814179
$ws.Cells.Item(12,12).Value = 1.181725805954586
$ws.Cells.Item(12,13).Value = 0.8727237325451784
$ws.Cells.Item(12,14).Value = 0.8727237325451784
$ws.Cells.Item(12,15).Value = 0.8170596412021528
$ws.Cells.Item(12,16).Value = 0.9089756568905916
$ws.Cells.Item(12,17).Value = 0.9089756568905916
$ws.Cells.Item(12,18).Value = 0.9271016190632981
$ws.Cells.Item(12,19).Value = 0.9271016190632981
$ws.Cells.Item(12,20).Value = 1.016693778860682
$ws.Cells.Item(13,3).Value = 1.029914601432402
$ws.Cells.Item(13,4).Value = 0.9730702420913527
$ws.Cells.Item(13,5).Value = 1.015245637183893
$ws.Cells.Item(13,6).Value = 1.012465677398484
$ws.Cells.Item(13,7).Value = 1.005285051518018
$ws.Cells.Item(13,8).Value = 0.9464665150187727
$ws.Cells.Item(13,9).Value = 0.9464665150187727
$ws.Cells.Item(13,10).Value = 1.029914601432402
$ws.Cells.Item(13,11).Value = 1.029914601432402
$ws.Cells.Item(13,12).Value = 1.005285051518018
$ws.Cells.Item(13,13).Value = 0.9758757832683953
$ws.Cells.Item(13,14).Value = 0.9758757832683953
$ws.Cells.Item(13,15).Value = 0.9749406028760478
$ws.Cells.Item(13,16).Value = 0.9938887226563974
$ws.Cells.Item(13,17).Value = 0.9938887226563975
$ws.Cells.Item(13,18).Value = 1.002895192350399
$ws.Cells.Item(13,19).Value = 1.002895192350399
$ws.Cells.Item(13,20).Value = 0.9970746207738205
$ws.Cells.Item(14,3).Value = 0.993244994999999
$ws.Cells.Item(14,4).Value = 0.9116258029473681
$ws.Cells.Item(14,5).Value = 1.164478286915789
$ws.Cells.Item(14,6).Value = 1.036995467957895
$ws.Cells.Item(14,7).Value = 1.054999769326315
$ws.Cells.Item(14,8).Value = 0.8700978064421067
$ws.Cells.Item(14,9).Value = 0.8700978064421067
$ws.Cells.Item(14,10).Value = 0.993244994999999
$ws.Cells.Item(14,11).Value = 0.993244994999999
$ws.Cells.Item(14,12).Value = 1.054999769326315
$ws.Cells.Item(14,13).Value = 0.9625487878842107
$ws.Cells.Item(14,14).Value = 0.9625487878842107
$ws.Cells.Item(14,15).Value = 0.9455744595719299
$ws.Cells.Item(14,16).Value = 0.9727808569228068
$ws.Cells.Item(14,17).Value = 0.9727808569228067
$ws.Cells.Item(14,18).Value = 0.9778968914421048
$ws.Cells.Item(14,19).Value = 0.9778968914421048
$ws.Cells.Item(14,20).Value = 1.005240354764912
$ws.Cells.Item(15,3).Value = 0.9802147653557584
$ws.Cells.Item(15,4).Value = 1.058601110701631
$ws.Cells.Item(15,5).Value = 0.921274419829571
$ws.Cells.Item(15,6).Value = 0.9736477442689017
$ws.Cells.Item(15,7).Value = 0.9709452632943285
$ws.Cells.Item(15,8).Value = 1.103561064059937
$ws.Cells.Item(15,9).Value = 1.103561064059937
$ws.Cells.Item(15,10).Value = 0.9802147653557584
$ws.Cells.Item(15,11).Value = 0.9802147653557584
$ws.Cells.Item(15,12).Value = 0.9709452632943285
$ws.Cells.Item(15,13).Value = 1.037253163677133
$ws.Cells.Item(15,14).Value = 1.037253163677133
$ws.Cells.Item(15,15).Value = 1.044369146018632
$ws.Cells.Item(15,16).Value = 1.018240364236675
$ws.Cells.Item(15,17).Value = 1.018240364236675
$ws.Cells.Item(15,18).Value = 1.008733964516445
$ws.Cells.Item(15,19).Value = 1.008733964516445
$ws.Cells.Item(15,20).Value = 1.001374061251688
$ws.Cells.Item(16,3).Value = 0.9618347199999993
$ws.Cells.Item(16,4).Value = 0.4918962299999994
$ws.Cells.Item(16,5).Value = 1.947381500000003
$ws.Cells.Item(16,6).Value = 1.212557700000001
$ws.Cells.Item(16,7).Value = 1.315735799999999
$ws.Cells.Item(16,8).Value = 0.2502035199999997
$ws.Cells.Item(16,9).Value = 0.2502035199999997
$ws.Cells.Item(16,10).Value = 0.9618347199999993
$ws.Cells.Item(16,11).Value = 0.9618347199999993
$ws.Cells.Item(16,12).Value = 1.315735799999999
$ws.Cells.Item(16,13).Value = 0.7829696599999996
$ws.Cells.Item(16,14).Value = 0.7829696599999996
$ws.Cells.Item(16,15).Value = 0.6859451833333328
$ws.Cells.Item(16,16).Value = 0.8425913466666661
$ws.Cells.Item(16,17).Value = 0.8425913466666661
$ws.Cells.Item(16,18).Value = 0.8724021899999994
$ws.Cells.Item(16,19).Value = 0.8724021899999994
$ws.Cells.Item(16,20).Value = 1.029934911666667
$ws.Cells.Item(17,3).Value = 0.99179523
$ws.Cells.Item(17,4).Value = 0.8102012999999999
$ws.Cells.Item(17,5).Value = 1.315746
$ws.Cells.Item(17,6).Value = 1.0863897
$ws.Cells.Item(17,7).Value = 1.1253174
$ws.Cells.Item(17,8).Value = 0.6733951500000001
$ws.Cells.Item(17,9).Value = 0.6733951500000001
$ws.Cells.Item(17,10).Value = 0.99179523
$ws.Cells.Item(17,11).Value = 0.99179523
$ws.Cells.Item(17,12).Value = 1.1253174
$ws.Cells.Item(17,13).Value = 0.8993562749999999
$ws.Cells.Item(17,14).Value = 0.8993562749999999
$ws.Cells.Item(17,15).Value = 0.8696379499999999
$ws.Cells.Item(17,16).Value = 0.93016926
$ws.Cells.Item(17,17).Value = 0.93016926
$ws.Cells.Item(17,18).Value = 0.9455757525
$ws.Cells.Item(17,19).Value = 0.9455757525
$ws.Cells.Item(17,20).Value = 1.00047413
$ws.Cells.Item(18,3).Value = 1.3679559
$ws.Cells.Item(18,4).Value = 0.7938064500000001
$ws.Cells.Item(18,5).Value = 0.9618347200000001
$ws.Cells.Item(18,6).Value = 1.1014624
$ws.Cells.Item(18,7).Value = 0.99179431
$ws.Cells.Item(18,8).Value = 0.5229061699999999
$ws.Cells.Item(18,9).Value = 0.5229061699999999
$ws.Cells.Item(18,10).Value = 1.3679559
$ws.Cells.Item(18,11).Value = 1.3679559
$ws.Cells.Item(18,12).Value = 0.99179431
$ws.Cells.Item(18,13).Value = 0.75735024
$ws.Cells.Item(18,14).Value = 0.75735024
$ws.Cells.Item(18,15).Value = 0.76950231
$ws.Cells.Item(18,16).Value = 0.9608854599999997
$ws.Cells.Item(18,17).Value = 0.9608854599999997
$ws.Cells.Item(18,18).Value = 1.06265307
$ws.Cells.Item(18,19).Value = 1.06265307
$ws.Cells.Item(18,20).Value = 0.9566266583333332
$ws.Cells.Item(19,3).Value = 0.8915424299999999
$ws.Cells.Item(19,4).Value = 1.1532739
$ws.Cells.Item(19,5).Value = 0.8283339199999999
$ws.Cells.Item(19,6).Value = 0.93542049
$ws.Cells.Item(19,7).Value = 0.9534773099999999
$ws.Cells.Item(19,8).Value = 1.2480078
$ws.Cells.Item(19,9).Value = 1.2480078
$ws.Cells.Item(19,10).Value = 0.8915424299999999
$ws.Cells.Item(19,11).Value = 0.8915424299999999
$ws.Cells.Item(19,12).Value = 0.9534773099999999
$ws.Cells.Item(19,13).Value = 1.100742555
$ws.Cells.Item(19,14).Value = 1.100742555
$ws.Cells.Item(19,15).Value = 1.118253003333333
$ws.Cells.Item(19,16).Value = 1.03100918
$ws.Cells.Item(19,17).Value = 1.03100918
$ws.Cells.Item(19,18).Value = 0.9961424924999999
$ws.Cells.Item(19,19).Value = 0.9961424924999999
$ws.Cells.Item(19,20).Value = 1.001675975
$ws.Cells.Item(20,3).Value = 1.177298675205479
$ws.Cells.Item(20,4).Value = 0.8021192230136984
$ws.Cells.Item(20,5).Value = 1.141209379178082
$ws.Cells.Item(20,6).Value = 1.09382152739726
$ws.Cells.Item(20,7).Value = 1.059468789863014
$ws.Cells.Item(20,8).Value = 0.5991860223287672
$ws.Cells.Item(20,9).Value = 0.5991860223287672
$ws.Cells.Item(20,10).Value = 1.177298675205479
$ws.Cells.Item(20,11).Value = 1.177298675205479
$ws.Cells.Item(20,12).Value = 1.059468789863014
$ws.Cells.Item(20,13).Value = 0.8293274060958904
$ws.Cells.Item(20,14).Value = 0.8293274060958904
$ws.Cells.Item(20,15).Value = 0.8202580117351598
$ws.Cells.Item(20,16).Value = 0.9453178291324201
$ws.Cells.Item(20,17).Value = 0.9453178291324199
$ws.Cells.Item(20,18).Value = 1.003313040650685
$ws.Cells.Item(20,19).Value = 1.003313040650685
$ws.Cells.Item(20,20).Value = 0.9788506028310501
$ws.Cells.Item(21,3).Value = 1.060644875789474
$ws.Cells.Item(21,4).Value = 0.9122313226315787
$ws.Cells.Item(21,5).Value = 1.079401156315789
$ws.Cells.Item(21,6).Value = 1.04188314
$ws.Cells.Item(21,7).Value = 1.034162258947369
$ws.Cells.Item(21,8).Value = 0.8205205531578945
$ws.Cells.Item(21,9).Value = 0.8205205531578945
$ws.Cells.Item(21,10).Value = 1.060644875789474
$ws.Cells.Item(21,11).Value = 1.060644875789474
$ws.Cells.Item(21,12).Value = 1.034162258947369
$ws.Cells.Item(21,13).Value = 0.9273414060526316
$ws.Cells.Item(21,14).Value = 0.9273414060526316
$ws.Cells.Item(21,15).Value = 0.9223047115789473
$ws.Cells.Item(21,16).Value = 0.9717758959649122
$ws.Cells.Item(21,17).Value = 0.9717758959649122
$ws.Cells.Item(21,18).Value = 0.9939931409210525
$ws.Cells.Item(21,19).Value = 0.9939931409210525
$ws.Cells.Item(21,20).Value = 0.991473884473684
$ws.Cells.Item(22,3).Value = 0.9584451305263157
$ws.Cells.Item(22,4).Value = 0.9144070715789473
$ws.Cells.Item(22,5).Value = 1.203648202631579
$ws.Cells.Item(22,6).Value = 1.032922123157895
$ws.Cells.Item(22,7).Value = 1.063571071578947
$ws.Cells.Item(22,8).Value = 0.9028521389473685
$ws.Cells.Item(22,9).Value = 0.9028521389473685
$ws.Cells.Item(22,10).Value = 0.9584451305263157
$ws.Cells.Item(22,11).Value = 0.9584451305263157
$ws.Cells.Item(22,12).Value = 1.063571071578947
$ws.Cells.Item(22,13).Value = 0.9832116052631579
$ws.Cells.Item(22,14).Value = 0.9832116052631579
$ws.Cells.Item(22,15).Value = 0.9602767607017544
$ws.Cells.Item(22,16).Value = 0.9749561136842105
$ws.Cells.Item(22,17).Value = 0.9749561136842105
$ws.Cells.Item(22,18).Value = 0.9708283678947368
$ws.Cells.Item(22,19).Value = 0.9708283678947368
$ws.Cells.Item(22,20).Value = 1.012640956403509
$ws.Cells.Item(23,3).Value = 1.035987231812905
$ws.Cells.Item(23,4).Value = 1.065926112515618
$ws.Cells.Item(23,5).Value = 0.8413008552341321
$ws.Cells.Item(23,6).Value = 0.9742366007178362
$ws.Cells.Item(23,7).Value = 0.948824807966586
$ws.Cells.Item(23,8).Value = 1.076880447286525
$ws.Cells.Item(23,9).Value = 1.076880447286525
$ws.Cells.Item(23,10).Value = 1.035987231812905
$ws.Cells.Item(23,11).Value = 1.035987231812905
$ws.Cells.Item(23,12).Value = 0.948824807966586
$ws.Cells.Item(23,13).Value = 1.012852627626555
$ws.Cells.Item(23,14).Value = 1.012852627626555
$ws.Cells.Item(23,15).Value = 1.030543789256243
$ws.Cells.Item(23,16).Value = 1.020564162355339
$ws.Cells.Item(23,17).Value = 1.020564162355339
$ws.Cells.Item(23,18).Value = 1.02441992971973
$ws.Cells.Item(23,19).Value = 1.02441992971973
$ws.Cells.Item(23,20).Value = 0.9905260092556002
$ws.Cells.Item(24,3).Value = 0.9652915770361605
$ws.Cells.Item(24,4).Value = 1.017857111946835
$ws.Cells.Item(24,5).Value = 1.011927935158048
$ws.Cells.Item(24,6).Value = 0.9896429273507429
$ws.Cells.Item(24,7).Value = 0.9996640585699197
$ws.Cells.Item(24,8).Value = 1.055007825826831
$ws.Cells.Item(24,9).Value = 1.055007825826831
$ws.Cells.Item(24,10).Value = 0.9652915770361605
$ws.Cells.Item(24,11).Value = 0.9652915770361605
$ws.Cells.Item(24,12).Value = 0.9996640585699197
$ws.Cells.Item(24,13).Value = 1.027335942198375
$ws.Cells.Item(24,14).Value = 1.027335942198375
$ws.Cells.Item(24,15).Value = 1.024176332114529
$ws.Cells.Item(24,16).Value = 1.006654487144304
$ws.Cells.Item(24,17).Value = 1.006654487144304
$ws.Cells.Item(24,18).Value = 0.996313759617268
$ws.Cells.Item(24,19).Value = 0.996313759617268
$ws.Cells.Item(24,20).Value = 1.006565239314756
$ws.Cells.Item(25,3).Value = 1.034299235804317
$ws.Cells.Item(25,4).Value = 1.065715635334983
$ws.Cells.Item(25,5).Value = 0.8459595168994624
$ws.Cells.Item(25,6).Value = 0.9736475715613275
$ws.Cells.Item(25,7).Value = 0.9486880292113041
$ws.Cells.Item(25,8).Value = 1.082247054752975
$ws.Cells.Item(25,9).Value = 1.082247054752975
$ws.Cells.Item(25,10).Value = 1.034299235804317
$ws.Cells.Item(25,11).Value = 1.034299235804317
$ws.Cells.Item(25,12).Value = 0.9486880292113041
$ws.Cells.Item(25,13).Value = 1.01546754198214
$ws.Cells.Item(25,14).Value = 1.01546754198214
$ws.Cells.Item(25,15).Value = 1.032216906433087
$ws.Cells.Item(25,16).Value = 1.021744773256199
$ws.Cells.Item(25,17).Value = 1.021744773256199
$ws.Cells.Item(25,18).Value = 1.024883388893228
$ws.Cells.Item(25,19).Value = 1.024883388893228
$ws.Cells.Item(25,20).Value = 0.9917595072607283
$ws.Cells.Item(26,3).Value = 0.9645573647888191
$ws.Cells.Item(26,4).Value = 1.018276349339149
$ws.Cells.Item(26,5).Value = 1.00942777778985
$ws.Cells.Item(26,6).Value = 0.9901223040686719
$ws.Cells.Item(26,7).Value = 1.000642851821473
$ws.Cells.Item(26,8).Value = 1.050096038616672
$ws.Cells.Item(26,9).Value = 1.050096038616672
$ws.Cells.Item(26,10).Value = 0.9645573647888191
$ws.Cells.Item(26,11).Value = 0.9645573647888191
$ws.Cells.Item(26,12).Value = 1.000642851821473
$ws.Cells.Item(26,13).Value = 1.025369445219072
$ws.Cells.Item(26,14).Value = 1.025369445219072
$ws.Cells.Item(26,15).Value = 1.023005079925764
$ws.Cells.Item(26,16).Value = 1.005098751742321
$ws.Cells.Item(26,17).Value = 1.005098751742321
$ws.Cells.Item(26,18).Value = 0.9949634050039458
$ws.Cells.Item(26,19).Value = 0.9949634050039458
$ws.Cells.Item(26,20).Value = 1.005520447737439
$ws.Cells.Item(27,3).Value = 1.032613787794711
$ws.Cells.Item(27,4).Value = 1.065504002113333
$ws.Cells.Item(27,5).Value = 0.8506176846503331
$ws.Cells.Item(27,6).Value = 0.9730590755090112
$ws.Cells.Item(27,7).Value = 0.9485509585659506
$ws.Cells.Item(27,8).Value = 1.08761119132621
$ws.Cells.Item(27,9).Value = 1.08761119132621
$ws.Cells.Item(27,10).Value = 1.032613787794711
$ws.Cells.Item(27,11).Value = 1.032613787794711
$ws.Cells.Item(27,12).Value = 0.9485509585659506
$ws.Cells.Item(27,13).Value = 1.01808107494608
$ws.Cells.Item(27,14).Value = 1.01808107494608
$ws.Cells.Item(27,15).Value = 1.033888717335165
$ws.Cells.Item(27,16).Value = 1.022925312562291
$ws.Cells.Item(27,17).Value = 1.02292531256229
$ws.Cells.Item(27,18).Value = 1.025347431370395
$ws.Cells.Item(27,19).Value = 1.025347431370395
$ws.Cells.Item(27,20).Value = 0.9929927833265916
$ws.Cells.Item(28,3).Value = 0.9638221813622511
$ws.Cells.Item(28,4).Value = 1.018695387507476
$ws.Cells.Item(28,5).Value = 1.006928559361844
$ws.Cells.Item(28,6).Value = 0.9906018303679908
$ws.Cells.Item(28,7).Value = 1.001622270446075
$ws.Cells.Item(28,8).Value = 1.045183547486109
$ws.Cells.Item(28,9).Value = 1.045183547486109
$ws.Cells.Item(28,10).Value = 0.9638221813622511
$ws.Cells.Item(28,11).Value = 0.9638221813622511
$ws.Cells.Item(28,12).Value = 1.001622270446075
$ws.Cells.Item(28,13).Value = 1.023402908966092
$ws.Cells.Item(28,14).Value = 1.023402908966092
$ws.Cells.Item(28,15).Value = 1.021833735146553
$ws.Cells.Item(28,16).Value = 1.003542666431478
$ws.Cells.Item(28,17).Value = 1.003542666431478
$ws.Cells.Item(28,18).Value = 0.9936125451641714
$ws.Cells.Item(28,19).Value = 0.9936125451641714
$ws.Cells.Item(28,20).Value = 1.004475629421957
$ws.Cells.Item(29,3).Value = 1.018200358817443
$ws.Cells.Item(29,4).Value = 0.9695990180661054
$ws.Cells.Item(29,5).Value = 1.029719632682923
$ws.Cells.Item(29,6).Value = 1.014727742239421
$ws.Cells.Item(29,7).Value = 1.013298665865324
$ws.Cells.Item(29,8).Value = 0.9354260382605486
$ws.Cells.Item(29,9).Value = 0.9354260382605486
$ws.Cells.Item(29,10).Value = 1.018200358817443
$ws.Cells.Item(29,11).Value = 1.018200358817443
$ws.Cells.Item(29,12).Value = 1.013298665865324
$ws.Cells.Item(29,13).Value = 0.9743623520629363
$ws.Cells.Item(29,14).Value = 0.9743623520629363
$ws.Cells.Item(29,15).Value = 0.9727745740639927
$ws.Cells.Item(29,16).Value = 0.9889750209811052
$ws.Cells.Item(29,17).Value = 0.9889750209811052
$ws.Cells.Item(29,18).Value = 0.9962813554401897
$ws.Cells.Item(29,19).Value = 0.9962813554401897
$ws.Cells.Item(29,20).Value = 0.9968285759886277
$ws.Cells.Item(30,3).Value = 1.032967370299977
$ws.Cells.Item(30,4).Value = 1.003742460472811
$ws.Cells.Item(30,5).Value = 0.9584463366816585
$ws.Cells.Item(30,6).Value = 0.9993618598241003
$ws.Cells.Item(30,7).Value = 0.9855324263618688
$ws.Cells.Item(30,8).Value = 0.9937901246824757
$ws.Cells.Item(30,9).Value = 0.9937901246824757
$ws.Cells.Item(30,10).Value = 1.032967370299977
$ws.Cells.Item(30,11).Value = 1.032967370299977
$ws.Cells.Item(30,12).Value = 0.9855324263618688
$ws.Cells.Item(30,13).Value = 0.9896612755221723
$ws.Cells.Item(30,14).Value = 0.9896612755221723
$ws.Cells.Item(30,15).Value = 0.9943550038390517
$ws.Cells.Item(30,16).Value = 1.004096640448107
$ws.Cells.Item(30,17).Value = 1.004096640448107
$ws.Cells.Item(30,18).Value = 1.011314322911075
$ws.Cells.Item(30,19).Value = 1.011314322911075
$ws.Cells.Item(30,20).Value = 0.9956400963871485
$ws.Cells.Item(31,3).Value = 0.9941579916311073
$ws.Cells.Item(31,4).Value = 1.099234551311955
$ws.Cells.Item(31,5).Value = 0.820806711715182
$ws.Cells.Item(31,6).Value = 0.9601642381183493
$ws.Cells.Item(31,7).Value = 0.9461750388630804
$ws.Cells.Item(31,8).Value = 1.134823384510616
$ws.Cells.Item(31,9).Value = 1.134823384510616
$ws.Cells.Item(31,10).Value = 0.9941579916311073
$ws.Cells.Item(31,11).Value = 0.9941579916311073
$ws.Cells.Item(31,12).Value = 0.9461750388630804
$ws.Cells.Item(31,13).Value = 1.040499211686848
$ws.Cells.Item(31,14).Value = 1.040499211686848
$ws.Cells.Item(31,15).Value = 1.06007765822855
$ws.Cells.Item(31,16).Value = 1.025052138334934
$ws.Cells.Item(31,17).Value = 1.025052138334934
$ws.Cells.Item(31,18).Value = 1.017328601658978
$ws.Cells.Item(31,19).Value = 1.017328601658978
$ws.Cells.Item(31,20).Value = 0.9925603193583816
